# Updates the cryptos list with the latest scraped price/volume snapshot.
# Column D (Price) and column E (Volume(1h)) are refreshed for most rows;
# rows 25/26 additionally swap which coin (Toncoin / Monero) sits at that
# rank, so B/C (name/link) change there too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. '1.009') are written with a leading apostrophe so they are stored
# as literal text, exactly like the other already-textual price strings.

$ws.Range('D2').Value = "28.092.59"
$ws.Range('E2').Value = "  +3.17%  "
$ws.Range('D3').Value = "1.881.40"
$ws.Range('E3').Value = "  +1.15%  "
$ws.Range('D4').Value = "'1.009"
$ws.Range('E4').Value = "  -1.20%  "
$ws.Range('D5').Value = "'313.87"
$ws.Range('E5').Value = "  +0.42%  "
$ws.Range('D6').Value = "'1.009"
$ws.Range('E6').Value = "  -0.99%  "
$ws.Range('D7').Value = "'0.4858"
$ws.Range('E7').Value = "  +1.44%  "
$ws.Range('D8').Value = "'0.3814"
$ws.Range('D9').Value = "'0.07368"
$ws.Range('E9').Value = "  +0.62%  "
$ws.Range('D10').Value = "'0.9440"
$ws.Range('E10').Value = "  +0.89%  "
$ws.Range('D11').Value = "'21.02"
$ws.Range('E11').Value = "  +3.40%  "
$ws.Range('D12').Value = "'0.07787"
$ws.Range('E12').Value = "  -0.99%  "
$ws.Range('D13').Value = "1.920.19"
$ws.Range('E13').Value = "  +2.87%  "
$ws.Range('D14').Value = "'5.538"
$ws.Range('E14').Value = "  +2.16%  "
$ws.Range('D15').Value = "'6.615"
$ws.Range('E15').Value = "  +1.20%  "
$ws.Range('D16').Value = "'91.76"
$ws.Range('E16').Value = "  +1.91%  "
$ws.Range('D17').Value = "'1.011"
$ws.Range('E17').Value = "  -1.10%  "
$ws.Range('D18').Value = "'0.000008877"
$ws.Range('E18').Value = "  +1.67%  "
$ws.Range('E19').Value = "  -0.88%  "
$ws.Range('D20').Value = "28.081.68"
$ws.Range('E20').Value = "  +3.01%  "
$ws.Range('D21').Value = "'14.91"
$ws.Range('E21').Value = "  +0.98%  "
$ws.Range('D22').Value = "'5.125"
$ws.Range('E22').Value = "  +0.37%  "
$ws.Range('D23').Value = "2.135.88"
$ws.Range('E23').Value = "  +1.57%  "
$ws.Range('D24').Value = "'10.99"
$ws.Range('E24').Value = "  +3.06%  "
$ws.Range('B25').Value = "Monero"
$ws.Range('C25').Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D25').Value = "'157.34"
$ws.Range('E25').Value = "  +2.31%  "
$ws.Range('B26').Value = "Toncoin"
$ws.Range('C26').Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('D26').Value = "'1.925"
$ws.Range('E26').Value = "  -1.48%  "
$ws.Range('D27').Value = "'18.57"
$ws.Range('E27').Value = "  +0.35%  "
$ws.Range('D28').Value = "'2.073"
$ws.Range('E28').Value = "  +3.65%  "
$ws.Range('D29').Value = "'115.82"
$ws.Range('E29').Value = "  +0.08%  "
$ws.Range('D30').Value = "'4.978"
$ws.Range('E30').Value = "  -0.20%  "
$ws.Range('D31').Value = "'0.08898"
$ws.Range('E31').Value = "  +0.12%  "
$ws.Range('D32').Value = "'3.333"
$ws.Range('E32').Value = "  -0.47%  "
$ws.Range('D33').Value = "'1.235"
$ws.Range('E33').Value = "  +4.04%  "
$ws.Range('D34').Value = "'0.7747"
$ws.Range('E34').Value = "  +4.65%  "
$ws.Range('D35').Value = "'4.657"
$ws.Range('E35').Value = "  +1.54%  "
$ws.Range('D36').Value = "'2.737"
$ws.Range('E36').Value = "  +2.06%  "
$ws.Range('D37').Value = "'1.129"
$ws.Range('E37').Value = "  +0.38%  "
$ws.Range('E38').Value = "  +0.80%  "
$ws.Range('D39').Value = "'0.5595"
$ws.Range('E39').Value = "  +4.90%  "
$ws.Range('E41').Value = "  +0.32%  "
$ws.Range('D42').Value = "'7.077"
$ws.Range('E42').Value = "  -0.50%  "
$ws.Range('D43').Value = "'8.555"
$ws.Range('E43').Value = "  +2.77%  "
$ws.Range('D44').Value = "'0.1525"
$ws.Range('E44').Value = "  -0.20%  "
$ws.Range('D45').Value = "'0.4899"
$ws.Range('E45').Value = "  +2.30%  "
$ws.Range('E46').Value = "  +0.95%  "
$ws.Range('D47').Value = "'105.75"
$ws.Range('E47').Value = "  +2.96%  "
$ws.Range('D48').Value = "'1.010"
$ws.Range('E48').Value = "  -0.99%  "
$ws.Range('D49').Value = "'1.674"
$ws.Range('E49').Value = "  +2.38%  "
$ws.Range('D50').Value = "'68.46"
$ws.Range('E50').Value = "  +3.02%  "
$ws.Range('D51').Value = "'0.06109"
$ws.Range('E51').Value = "  +0.54%  "
